# Update underlying "PO List" source data. The "Assignment" sheet's E10 cell
# holds a dynamic array FILTER() formula that picks the name with the highest
# 'ppvc_next' (column V) value among qualifying rows (W=1 and G=1); updating
# the source data and recalculating lets Excel re-derive the new cached
# result ("Woo Kwan Wye") itself, rather than overwriting the formula cell
# directly (which would destroy the formula).

$wb = $excel.ActiveWorkbook
$wsPOList = $wb.Worksheets.Item("PO List")

# Row 3 (Catherine Wong) - 10storey_next
$wsPOList.Range("K3").Value = 10

# Row 4 (Charmaine Fang) - ppvc_next
$wsPOList.Range("V4").Value = 3

# Row 6 (Richard Lim) - 10storey_next
$wsPOList.Range("K6").Value = 13

# Row 9 (Michael Sien) - multiple fields updated
$wsPOList.Range("H9").Value = 2
$wsPOList.Range("I9").Value = 2
$wsPOList.Range("J9").Value = 44900
$wsPOList.Range("K9").Value = 2
$wsPOList.Range("P9").Value = 1
$wsPOList.Range("Q9").Value = 1
$wsPOList.Range("R9").Value = 44900
$wsPOList.Range("S9").Value = 2
$wsPOList.Range("T9").Value = 2
$wsPOList.Range("U9").Value = 44900
$wsPOList.Range("V9").Value = 1

# Row 10 (Lim Meng May) - 10storey_next
$wsPOList.Range("K10").Value = 8

# Row 11 (Woo Kwan Wye) - 10storey_next, ppvc_next
$wsPOList.Range("K11").Value = 15
$wsPOList.Range("V11").Value = 6

# Row 12 (Brian Phua) - multiple fields updated
$wsPOList.Range("H12").Value = 8
$wsPOList.Range("I12").Value = 8
$wsPOList.Range("J12").Value = 44873
$wsPOList.Range("K12").Value = 4
$wsPOList.Range("P12").Value = 3
$wsPOList.Range("Q12").Value = 3
$wsPOList.Range("R12").Value = 44433
$wsPOList.Range("S12").Value = 7
$wsPOList.Range("T12").Value = 7
$wsPOList.Range("U12").Value = 44873
$wsPOList.Range("V12").Value = 2

# Row 15 (Lim Leng Boon) - 10storey_next
$wsPOList.Range("K15").Value = 9

# Row 16 (Tan Chong Lin) - 10storey_next, ppvc_next
$wsPOList.Range("K16").Value = 12
$wsPOList.Range("V16").Value = 5

# Row 17 (Victor Tay) - 10storey_next
$wsPOList.Range("K17").Value = 11

# Row 19 (Kang Mi) - 10storey_next, ppvc_next
$wsPOList.Range("K19").Value = 7
$wsPOList.Range("V19").Value = 4

# Row 24 (Zhao Zhe) - 10storey_next
$wsPOList.Range("K24").Value = 5

# Row 26 (Willie Chai) - 10storey_next
$wsPOList.Range("K26").Value = 13

# Row 29 (Chang Heng Choy) - 10storey_next
$wsPOList.Range("K29").Value = 6

# Recalculate the whole workbook so dependent formulas (e.g. Assignment!E10)
# pick up the refreshed cached values.
$excel.CalculateFullRebuild()
